$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RelivePos" column (E) values changed from "186,6.89,88" to "186,0,88"
# for every data row on the sheet (rows 2-4).
$ws.Range("E2").Value = "186,0,88"
$ws.Range("E3").Value = "186,0,88"
$ws.Range("E4").Value = "186,0,88"

# Move / update the active selection from E4 to F7.
$ws.Range("F7").Select() | Out-Null
